$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.963.34"
$ws.Range("E2").Value = "  -1.83%  "

$ws.Range("D3").Value = "'1.868.61"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").Value = "'312.30"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "'0.4993"
$ws.Range("E7").Value = "  -2.36%  "

$ws.Range("D8").Value = "'0.3820"
$ws.Range("E8").Value = "  -3.59%  "

$ws.Range("D9").Value = "'0.08990"
$ws.Range("E9").Value = "  -7.51%  "

$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").Value = "'41.53"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").Value = "'6.381"
$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("D13").Value = "'20.69"
$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("D14").Value = "'1.861.32"
$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "'7.234"
$ws.Range("E15").Value = "  -2.58%  "

$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "'0.00001102"
$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("D18").Value = "'91.02"
$ws.Range("E18").Value = "  -3.37%  "

$ws.Range("D19").Value = "'0.06665"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").Value = "'18.04"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").Value = "'6.119"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "'28.004.34"
$ws.Range("E23").Value = "  -2.09%  "

$ws.Range("D24").Value = "'11.54"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").Value = "'2.281"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.390"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.086.50"
$ws.Range("E27").Value = "  -2.34%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.514"
$ws.Range("E28").Value = "  -5.83%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'158.50"
$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.73"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'126.21"
$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1060"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.059"
$ws.Range("E33").Value = "  -4.34%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.580"
$ws.Range("E34").Value = "  -1.61%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.579"
$ws.Range("E35").Value = "  -1.73%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.412"
$ws.Range("E36").Value = "  -3.47%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06547"
$ws.Range("E37").Value = "  -2.49%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02424"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").Value = "'0.2196"
$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.291"
$ws.Range("E40").Value = "  +9.12%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.201"
$ws.Range("E41").Value = "  -4.49%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.61"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6374"
$ws.Range("E43").Value = "  -0.67%  "

$ws.Range("B44").Value = "InternetComputer(DFINITY)"
$ws.Range("C44").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D44").Value = "'4.901"
$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.23"
$ws.Range("E46").Value = "  -2.14%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6010"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "'3.669"
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.274"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").Value = "'1.997"
$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = "  +1.97%  "
